# "Generate Report for Handback" - mark the two e2e files as handed back
# and record the handback datetime/file for each language sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: the localization status text changes for both languages ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- Per-language sheets: zh-cn and de-de ---
$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/005232eac6de5961dc9546f12d1cc8258348e6e1/e2e/"

$langs = @(
    @{ Name = "zh-cn"; HandbackDate = "2016-08-23 14:26:48" },
    @{ Name = "de-de"; HandbackDate = "2016-08-23 14:26:56" }
)

$files = @(
    @{ Row = 2; Md = "42bdd63f-3fe6-402d-92c0-9d5c4b45fa99.md" },
    @{ Row = 3; Md = "cbcdcd7e-44cf-41bd-bb70-60f85ff85b5e.md" }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Name)

    foreach ($file in $files) {
        $row = $file.Row

        # Status column (C)
        $ws.Range("C$row").Value = $newStatus

        # Latest Handoff File (G) content -> also becomes Latest Handback File (J)
        $handoffFile = $ws.Range("G$row").Value2

        # Latest Target File (I) = the source .md file, as a hyperlink
        $ws.Range("I$row").Value = $file.Md
        $ws.Hyperlinks.Add($ws.Range("I$row"), ($baseUrl + $file.Md), "", "", $file.Md) | Out-Null

        # Latest Handback File (J) = same xlf file recorded as handoff
        $ws.Range("J$row").Value = $handoffFile

        # Latest Handback DateTime (K)
        $ws.Range("K$row").Value = $lang.HandbackDate
    }
}
